# ResourceMonitor v1.1.1 - log threshold condition added
# Updates the "기능" (Feature) sheet:
#  - Win32_PerfRawData_PerfProc_Process -> Win32_PerfFormattedData_PerfProc_Process (E3:E5)
#  - resets the stray large/dark font on E5 back to the sheet's normal formatting
#  - normalises the B11:B13 group to the same left/vcenter alignment used by every
#    other group in column B (it was previously centered)
#  - appends three new rows (24-26) describing the new "log setting" feature
#    (log interval + CPU/memory threshold), with matching merges & alignment

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) E3/E4/E5: WMI class rename (Raw -> Formatted perf data)
# ---------------------------------------------------------------------------
$ws.Range("E3").Value = "Win32_PerfFormattedData_PerfProc_Process"
$ws.Range("E4").Value = "Win32_PerfFormattedData_PerfProc_Process"
$ws.Range("E5").Value = "Win32_PerfFormattedData_PerfProc_Process"

# E5 previously carried a one-off font override (12pt / dark grey). Bring it
# back in line with the rest of the row (plain/default format) by copying the
# plain formatting of its neighbour F5 onto it.
$ws.Range("F5").Copy()
$ws.Range("E5").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Row 5 also had a custom row height tied to that font override - let it
# return to the sheet's standard height.
$ws.Rows.Item(5).AutoFit()

# ---------------------------------------------------------------------------
# 2) B11:B13 used a centred style left over from an older layout; every other
#    feature group (B2, B6, B14, B18, ...) uses left/vcenter - match it.
# ---------------------------------------------------------------------------
$ws.Range("B11:B13").HorizontalAlignment = -4131
$ws.Range("B11:B13").VerticalAlignment = -4108

# ---------------------------------------------------------------------------
# 3) New rows 24-26: "로그 설정" (log setting) feature block
# ---------------------------------------------------------------------------
$ws.Range("A24").Value = "로그 설정"
$ws.Range("B24").Value = "로그 주기 설정"
$ws.Range("D24").Value = "ms"
$ws.Range("E24").Value = "CDlgSetLogInterval"

$ws.Range("B25").Value = "로그 조건 설정"
$ws.Range("C25").Value = "CPU Threshold 설정 "
$ws.Range("D25").Value = "%"
$ws.Range("E25").Value = "DlgSetLogThreshold"

$ws.Range("C26").Value = "메모리 Threshold 설정 "
$ws.Range("D26").Value = "MB"
$ws.Range("E26").Value = "DlgSetLogThreshold"

# Merge the grouped cells first ...
$ws.Range("A24:A26").Merge()
$ws.Range("B24:C24").Merge()
$ws.Range("B25:B26").Merge()

# ... then apply the same alignment conventions used throughout the sheet:
#  - column A group label: left / vcenter
$ws.Range("A24:A26").HorizontalAlignment = -4131
$ws.Range("A24:A26").VerticalAlignment = -4108

#  - B24:C24 feature name header above the sub rows: left / top
$ws.Range("B24:C24").HorizontalAlignment = -4131
$ws.Range("B24:C24").VerticalAlignment = -4160

#  - B25:B26 sub-feature label: left / vcenter
$ws.Range("B25:B26").HorizontalAlignment = -4131
$ws.Range("B25:B26").VerticalAlignment = -4108

#  - C25/C26 detail labels: left / vcenter
$ws.Range("C25:C26").HorizontalAlignment = -4131
$ws.Range("C25:C26").VerticalAlignment = -4108

#  - D24/D25 unit cells: vcenter only (matches D2, D3, D8, ... pattern)
$ws.Range("D24:D25").VerticalAlignment = -4108

# ---------------------------------------------------------------------------
# 4) Selection / viewport left the way the author saved the sheet
# ---------------------------------------------------------------------------
$ws.Range("B25:B26").Select()
